$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend column N with the same formatting as column M (borders / number
# formats / fonts for rows 3-14), then fill in the new 2020 column values.
$ws.Range("M3:M14").Copy()
$ws.Range("N3:N14").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# N3 (the blank cell in the bottom-border row) uses the bold/no-theme-color
# style from the row-14 label cells rather than the plain M3 style.
$ws.Range("A14").Copy()
$ws.Range("N3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New column header (year 2020)
$ws.Range("N4").Value = 2020

# New column N data values
$ws.Range("N5").Value = 68.5
$ws.Range("N6").Value = 106.7
$ws.Range("N7").Value = 53.2
$ws.Range("N8").Value = 49.6
$ws.Range("N9").Value = 108.9
$ws.Range("N10").Value = 107.8
$ws.Range("N11").Value = 155.7
$ws.Range("N12").Value = 25.9
$ws.Range("N13").Value = 103.5
$ws.Range("N14").Value = 11

# Revised values in existing column L / M cells
$ws.Range("M5").Value = 68.4
$ws.Range("M6").Value = 108.2
$ws.Range("M7").Value = 51.7
$ws.Range("M8").Value = 97.7
$ws.Range("L9").Value = 105.6
$ws.Range("M9").Value = 106.7
$ws.Range("M10").Value = 124.2
$ws.Range("M11").Value = 138.8
$ws.Range("L12").Value = 27.1
$ws.Range("M12").Value = 33.9
$ws.Range("M13").Value = 96
$ws.Range("M14").Value = 7.7

# The source workbook also carries a page setup (print settings) that
# wasn't present before the edit.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
